# Refactor to agentic architecture, add executive dashboard, UI integration,
# and audit-safe pipeline.
#
# - Rename "file" header to "file_name"
# - Add 5 new audit/quality columns (I:M)
# - Replace the detected invoice number with a generic placeholder
# - Convert processed_utc from a literal text timestamp to a real Excel
#   date/time serial value (custom number format)
# - Append 3 more rows of (low-confidence / on-hold) review data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------
$ws.Cells.Item(1, 1).Value = "file_name"
$ws.Cells.Item(1, 9).Value  = "ink_density_score"
$ws.Cells.Item(1, 10).Value = "avg_ocr_confidence"
$ws.Cells.Item(1, 11).Value = "audit_remarks"
$ws.Cells.Item(1, 12).Value = "decision_stage"
$ws.Cells.Item(1, 13).Value = "data_quality_flag"

# ---- Row 2 : existing record, updated ------------------------------
$ws.Cells.Item(2, 2).Value = "AUTO-DETECTED"

# processed_utc -> real datetime serial, custom format
$ws.Range("H2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("H2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 8).Value = 46009.45200974537

$ws.Cells.Item(2, 9).Value = ""
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 12).Value = ""
$ws.Cells.Item(2, 13).Value = ""

# ---- Row 3 : new low-confidence record ------------------------------
$ws.Cells.Item(3, 1).Value = "sample1.PNG"
$ws.Cells.Item(3, 2).Value = ""
$ws.Cells.Item(3, 3).Value = 642
$ws.Cells.Item(3, 4).Value = "INR"
$ws.Cells.Item(3, 5).Value = $true
$ws.Cells.Item(3, 6).Value = 0.4
$ws.Cells.Item(3, 7).Value = "ON_HOLD"

$ws.Range("H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 8).Value = 46015.3009153588

$ws.Cells.Item(3, 9).Value = 0.1215
$ws.Cells.Item(3, 10).Value = 0.174
$ws.Cells.Item(3, 11).Value = "Low OCR confidence"
$ws.Cells.Item(3, 12).Value = "REQUIRES_MANUAL_REVIEW"
$ws.Cells.Item(3, 13).Value = "LOW_CONFIDENCE"

# ---- Row 4 : new low-confidence record ------------------------------
$ws.Cells.Item(4, 1).Value = "sample1.PNG"
$ws.Cells.Item(4, 2).Value = ""
$ws.Cells.Item(4, 3).Value = 642
$ws.Cells.Item(4, 4).Value = "INR"
$ws.Cells.Item(4, 5).Value = $true
$ws.Cells.Item(4, 6).Value = 0.4
$ws.Cells.Item(4, 7).Value = "ON_HOLD"

$ws.Range("H4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 8).Value = 46015.30579377315

$ws.Cells.Item(4, 9).Value = 0.1215
$ws.Cells.Item(4, 10).Value = 0.174
$ws.Cells.Item(4, 11).Value = "Low OCR confidence"
$ws.Cells.Item(4, 12).Value = "REQUIRES_MANUAL_REVIEW"
$ws.Cells.Item(4, 13).Value = "LOW_CONFIDENCE"

# ---- Row 5 : new low-confidence record ------------------------------
$ws.Cells.Item(5, 1).Value = "sample1.PNG"
$ws.Cells.Item(5, 2).Value = ""
$ws.Cells.Item(5, 3).Value = 642
$ws.Cells.Item(5, 4).Value = "INR"
$ws.Cells.Item(5, 5).Value = $true
$ws.Cells.Item(5, 6).Value = 0.4
$ws.Cells.Item(5, 7).Value = "ON_HOLD"

$ws.Range("H5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 8).Value = 46015.30790904717

$ws.Cells.Item(5, 9).Value = 0.1215
$ws.Cells.Item(5, 10).Value = 0.174
$ws.Cells.Item(5, 11).Value = "Low OCR confidence"
$ws.Cells.Item(5, 12).Value = "REQUIRES_MANUAL_REVIEW"
$ws.Cells.Item(5, 13).Value = "LOW_CONFIDENCE"
